$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The table "Tabela1" holds the daily COVID-19 stats; the update appends
# one more day's row (2020-05-26, date serial 43976) to the bottom of the
# table, growing it from A1:J75 to A1:J76.
$lo = $ws.ListObjects.Item("Tabela1")
$newRow = $lo.ListRows.Add()
$rng = $newRow.Range

# Clone the formatting of the previous data row (date format, fonts,
# alignment, …) onto the new one before filling in the values, so the new
# row renders exactly like the rest of the table instead of picking up
# generic default formatting.
$ws.Range("A75:J75").Copy()
$rng.PasteSpecial(-4122)
$excel.CutCopyMode = 0

$rng.Item(1, 1).Value = 43976
$rng.Item(1, 2).Value = 75770
$rng.Item(1, 3).Value = 754
$rng.Item(1, 4).Value = 1469
$rng.Item(1, 5).Value = 0
$rng.Item(1, 6).Value = 9
$rng.Item(1, 7).Value = 2
$rng.Item(1, 8).Value = 6
$rng.Item(1, 9).Value = 108
$rng.Item(1, 10).Value = 1

# Match the author's final selection (the newly added row).
$ws.Range("A76:J76").Select()
